# Roll the IH futures-contract schedule forward by one quarter and clear
# the stray trailing (empty) row that was left over from the previous sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IH")

$ws.Range("B2").Value = "IH2306.CFFEX,IH2309.CFFEX"
$ws.Range("B3").Value = "IH2309.CFFEX,IH2312.CFFEX"
$ws.Range("B4").Value = "IH2312.CFFEX,IH2403.CFFEX"
$ws.Range("B5").Value = "IH2403.CFFEX,IH2406.CFFEX"
$ws.Range("B6").Value = "IH2406.CFFEX,IH2409.CFFEX"

$ws.Range("A7:B7").ClearContents()

$ws.Columns.Item(2).AutoFit() | Out-Null
